# ScopusCitedByCore.xlsx update:
#  - fill in raw-input columns H:O for rows 74-88 (previously blank -> zeros)
#  - append new daily rows 89-112 (dates + formulas + raw inputs)
#  - append summary rows 115-117 (Average / Min / Max) with new styling
#  - nudge the saved window position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Raw input data (columns H:O) for rows 74-112.
#    Columns B:F on rows 74-88 already hold formulas (=L,=M,=I,=N-M,=O-M) that
#    will simply pick up the new inputs once H:O are populated.
# ---------------------------------------------------------------------------
$rawData = @(
    @(74,  55301184, 1330, 19846,  69729,  55300994, 1140,  19656,  69539),
    @(75,  55300994, 1174, 290700, 212989, 55302145, 2325,  291851, 214140),
    @(76,  55302145, 2199, 79037,  349720, 55312630, 12684, 89522,  360205),
    @(77,  55312630, 2837, 18835,  68134,  55317727, 7934,  23932,  73231),
    @(78,  55317727, 206,  58672,  37083,  55318636, 1115,  59581,  37992),
    @(79,  55318636, 719,  118766, 170987, 55336470, 18553, 136600, 188821),
    @(80,  55336470, 116,  285476, 8698,   55336596, 242,   285602, 8824),
    @(81,  55336596, 210,  75661,  42799,  55337818, 1432,  76883,  44021),
    @(82,  55337818, 2765, 169373, 387966, 55360779, 25726, 192334, 410927),
    @(83,  55360779, 5407, 184537, 245715, 55373139, 17767, 196897, 258075),
    @(84,  55373139, 1812, 21998,  17049,  55371743, 416,   20602,  15653),
    @(85,  55371743, 0,    110880, 170776, 55384014, 12271, 123151, 183047),
    @(86,  55384014, 852,  147064, 683899, 55405454, 22292, 168504, 705339),
    @(87,  55405454, 57,   57,     881,    55405397, 0,     0,      824),
    @(88,  55405397, 0,    0,      93,     55405399, 2,     2,      95),
    @(89,  55405399, 2,    185767, 173536, 55416799, 11402, 197167, 184936),
    @(90,  55416799, 1953, 295799, 597101, 55435486, 20640, 314486, 615788),
    @(91,  55435486, 5906, 116042, 191505, 55449703, 20123, 130259, 205722),
    @(92,  55449703, 1153, 192887, 166442, 55461346, 12796, 204530, 178085),
    @(93,  55461346, 4355, 241548, 77145,  55459172, 2181,  239374, 74971),
    @(94,  55459172, 2214, 22275,  234467, 55469211, 12253, 32314,  244506),
    @(95,  55469211, 4615, 161446, 273227, 55480394, 15798, 172629, 284410),
    @(96,  55480394, 5425, 160402, 174243, 55490490, 15521, 170498, 184339),
    @(97,  55490490, 3123, 160739, 335971, 55513629, 26262, 183878, 359110),
    @(98,  55513629, 42,   231506, 248800, 55521256, 7669,  239133, 256427),
    @(99,  55521256, 4048, 221861, 469115, 55532580, 15372, 233185, 480439),
    @(100, 55532580, 2890, 156363, 24530,  55540007, 10317, 163790, 31957),
    @(101, 55540007, 8633, 24261,  15358,  55531374, 0,     15628,  6725),
    @(102, 55531374, 2627, 174661, 873637, 55544627, 15880, 187914, 886890),
    @(103, 55544627, 2113, 288429, 959240, 55568540, 26026, 312342, 983153),
    @(104, 55568540, 2189, 129695, 149058, 55591619, 25268, 152774, 172137),
    @(105, 55591619, 329,  65786,  10718,  55591992, 702,   66159,  11091),
    @(106, 55591992, 1857, 143010, 313036, 55616660, 26525, 167678, 337704),
    @(107, 55616660, 2752, 210927, 272796, 55635261, 21353, 229528, 291397),
    @(108, 55635261, 1963, 77624,  570034, 55682268, 48970, 124631, 617041),
    @(109, 55682268, 5,    22006,  315969, 55706099, 23836, 45837,  339800),
    @(110, 55706099, 0,    0,      0,      55706099, 0,     0,      0),
    @(111, 55706099, 0,    0,      0,      55706099, 0,     0,      0),
    @(112, 55706099, 21675,164112, 963901, 55743891, 59467, 201904, 1001693)
)

foreach ($entry in $rawData) {
    $r = $entry[0]
    $arr = New-Object 'object[,]' 1,8
    for ($i = 0; $i -lt 8; $i++) { $arr[0,$i] = $entry[$i + 1] }
    $rng = "H" + $r + ":O" + $r
    $ws.Range($rng).Value = $arr
}

# ---------------------------------------------------------------------------
# 2. New rows 89-112: date in column A, shared-style formulas in B:F.
# ---------------------------------------------------------------------------
for ($r = 89; $r -le 112; $r++) {
    $dateSerial = 41671 + ($r - 89)
    $ws.Range("A" + $r).Value = $dateSerial
    $ws.Range("B" + $r).Formula = "=L" + $r
    $ws.Range("C" + $r).Formula = "=M" + $r
    $ws.Range("D" + $r).Formula = "=I" + $r
    $ws.Range("E" + $r).Formula = "=N" + $r + "-M" + $r
    $ws.Range("F" + $r).Formula = "=O" + $r + "-M" + $r
}
# Match the date-style ("s=2") already used by rows 58-88 in column A.
$ws.Range("A89:A112").NumberFormat = "[$-409]d\-mmm;@"

# ---------------------------------------------------------------------------
# 3. Summary rows: Average / Min / Max over C3:F112.
# ---------------------------------------------------------------------------
$ws.Range("A115").Value = "Average"
$ws.Range("A116").Value = "Min"
$ws.Range("A117").Value = "Max"

$ws.Range("C115").Formula = "=AVERAGE(C3:C112)"
$ws.Range("D115:F115").Formula = "=AVERAGE(D3:D112)"

$ws.Range("C116").Formula = "=MIN(C3:C112)"
$ws.Range("D116:F116").Formula = "=MIN(D3:D112)"

$ws.Range("C117").Formula = "=MAX(C3:C112)"
$ws.Range("D117:F117").Formula = "=MAX(D3:D112)"

# New font (black RGB instead of the automatic theme color) across the block.
$ws.Range("A115:F117").Font.Color = 0
# Average row's numeric cells get an integer number format; Min/Max keep General.
$ws.Range("C115:F115").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4. Cosmetic view state (best effort - window chrome position / scroll).
# ---------------------------------------------------------------------------
$ws.Range("G119").Select()
try { $excel.ActiveWindow.ScrollRow = 93 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
try { $excel.ActiveWindow.Left = 5800 } catch {}
